# Daily attendance processing - 2025-12-29 17:31:41
# Swap the "Recorded By" display order from "System, <email>" to
# "<email>, System" for the specific session rows that were re-processed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2,3,4,5,6,7,16,17,21,22,23,36,37,41,42,43,56,57,61,62,63,76,77,81,82,83,84,85,86,95,96,100,101,102,103,104,105,114,115,119,120,121,122,123,124,133,134,138,139,140,141,142,143,152,153,157,158,159,160,161,162,171,172,176,177,178,191,192,196,197,198,211,212,216,217,218,231,232)

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
